# Naan mudhalvan.pptx - "Add files via upload" commit replay
#
# Changes applied:
#  1) Date placeholder fields (notes master, slide layout, slide master)
#     9/22/2023 -> 11/17/2023
#  2) Slide 1, shape "object 3": split the "Task - 1" run into "Task - " + "1"
#     (same formatting) so the "1" can carry its own run.
#  3) Slide 19, shape "object 4": shrink the textbox height slightly, fix the
#     run language to en-US, and swap the GitHub submission URL for the new
#     repository link.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Date fields: 9/22/2023 -> 11/17/2023
# ---------------------------------------------------------------------------

# Notes master -> "Date Placeholder 2"
$nm = $p.NotesMaster
$nmDate = $nm.Shapes.Item(2)
$nmDate.TextFrame.TextRange.Text = "11/17/2023"

# Slide master -> "Holder 5"
$sm = $p.SlideMaster
$smDate = $sm.Shapes.Item(4)
$smDate.TextFrame.TextRange.Text = "11/17/2023"

# Slide layout (CustomLayout 1) -> "Date 3"
$cl = $sm.CustomLayouts.Item(1)
$clDate = $cl.Shapes.Item(3)
$clDate.TextFrame.TextRange.Text = "11/17/2023"

# ---------------------------------------------------------------------------
# 2) Slide 1 - "Task - 1" run split
# ---------------------------------------------------------------------------

$s1 = $p.Slides.Item(1)
$taskShape = $s1.Shapes.Item(2)
$taskText = $taskShape.TextFrame.TextRange
$taskPara = $taskText.Paragraphs(2)

# Re-apply identical character formatting to the trailing "1" so PowerPoint
# breaks it into its own run, matching the "Task - " / "1" run split.
$taskDigit = $taskPara.Characters(8, 1)
$taskDigit.Font.Bold = -1
$taskDigit.Font.Size = 24
$taskDigit.Font.Name = "CFJCTS+PublicSans-Bold"
$taskDigit.Font.Color.RGB = 6895138

# ---------------------------------------------------------------------------
# 3) Slide 19 - GitHub submission link textbox
# ---------------------------------------------------------------------------

$s19 = $p.Slides.Item(19)
$linkShape = $s19.Shapes.Item(3)

# Shrink the textbox height to match the new (shorter) wrapped text.
$linkShape.Height = 30.974566929133857

$linkText = $linkShape.TextFrame.TextRange
$linkRun = $linkText.Runs(1)
$linkRun.Text = "https://github.com/mahendran121/NM-SPCET-CSE-GROUP06"
$linkRun.LanguageID = "en-US"
